# Update "想去人数" (F column) figures for the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

# Row -> new value for column F, for sheet "展览" (sheet1)
$exhibitionUpdates = @{
    3  = 346
    5  = 1767
    6  = 91
    7  = 2225
    9  = 286
    11 = 5034
    12 = 370
    14 = 312
    16 = 35
    17 = 194
    18 = 388
    21 = 4029
    22 = 725
    23 = 711
    27 = 128
    31 = 591
    34 = 1031
    36 = 2620
    38 = 34
}

# Row -> new value for column F, for sheet "全部类型" (sheet4)
$allTypesUpdates = @{
    3  = 346
    5  = 1767
    6  = 91
    7  = 2225
    9  = 286
    11 = 5034
    12 = 370
    14 = 312
    16 = 35
    17 = 194
    18 = 388
    21 = 4029
    22 = 725
    23 = 711
    27 = 128
    31 = 591
    35 = 1031
    37 = 2620
    39 = 34
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Cells.Item($row, 6).Value = $exhibitionUpdates[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Cells.Item($row, 6).Value = $allTypesUpdates[$row]
}
